$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-10-18 Saturday" "2025-10-19 Sunday"

Replace-Text "45÷8=" "94÷4="
Replace-Text "14÷8=" "84÷4="
Replace-Text "19÷2=" "90÷5="
Replace-Text "59÷6=" "78÷5="
Replace-Text "75÷7=" "67÷5="
Replace-Text "79÷2=" "33÷8="
Replace-Text "76÷6=" "10÷7="
Replace-Text "22÷6=" "49÷7="
Replace-Text "99÷9=" "92÷8="
Replace-Text "65÷8=" "56÷7="
Replace-Text "47÷8=" "88÷5="
Replace-Text "18÷4=" "80÷4="
Replace-Text "74÷6=" "75÷4="
Replace-Text "21÷6=" "42÷6="
Replace-Text "26÷9=" "65÷4="
Replace-Text "62÷5=" "36÷7="
Replace-Text "42÷3=" "40÷7="
Replace-Text "14÷9=" "11÷2="
Replace-Text "17÷4=" "41÷5="
Replace-Text "70÷9=" "38÷9="
Replace-Text "18÷6=" "71÷8="
Replace-Text "93÷4=" "13÷8="
Replace-Text "60÷3=" "74÷5="
Replace-Text "65÷9=" "15÷7="
Replace-Text "82÷4=" "80÷7="
